$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.329.31'
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = '  +3.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.925.20'
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = '  +2.35%  '

$ws.Range("E4").Value = '  -0.88%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.35'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("E7").Value = '  +1.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3868'
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = '  +2.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07403'
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9452'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.94'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07840'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.935.07'
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = '  +3.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.541'
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = '  +1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.648'
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = '  +0.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.84'
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = '  +0.98%  '

$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008883'
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  -0.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.346.32'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = '  +3.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.98'
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.169'
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.133.72'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +1.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.97'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  +2.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.936'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = '  -1.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.21'
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = '  +1.38%  '

$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.106'
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = '  +3.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.92'
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.013'
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = '  -0.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08919'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.366'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.260'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  +2.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7825'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  +4.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.720'
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = '  +2.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.758'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = '  +1.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02063'
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = '  -0.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.116'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05372'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  +1.09%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5574'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = '  +3.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.039'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = '  +0.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.080'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1540'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  +0.15%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.571'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = '  +1.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4917'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = '  +1.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.72'
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.10'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.676'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  +0.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.48'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  +3.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06146'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  +0.40%  '
